$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "Bapco closure (Y/N)" column (H) for all comment rows with "Y"
$ws.Range("H2").Value = "Y"
$ws.Range("H3").Value = "Y"
$ws.Range("H4").Value = "Y"
$ws.Range("H5").Value = "Y"
$ws.Range("H6").Value = "Y"

# Row 5's "Comment included (Y/N)" (G) was N.A. but is now marked Y
$ws.Range("G5").Value = "Y"

# Restore active cell selection
$ws.Range("G5").Select()
